# Updated symbol list on Wed Jan 18 23:55:03 UTC 2023 with GitHub Actions
# Refresh Price (col D) and Volume(1h) (col E) figures for the crypto table.
# Values are assigned with a leading apostrophe (and Style reset to "Normal"
# afterwards) so they stay plain text cells, matching the original sheet's
# inline-string "Price"/"Volume(1h)" columns instead of being coerced to
# numeric/percentage values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'287.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-4.39%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'30.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-4.04%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.895"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.34%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07100"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-9.90%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.791"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-13.65%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.644"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-1.97%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.779"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-1.54%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.8962"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-3.18%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1640"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-5.69%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07525"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-5.42%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07985"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-7.86%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.02999"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-3.20%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09993"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.17%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001499"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-1.20%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005640"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-4.30%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D18").Value = "'3.469"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.32%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-6.02%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'0.00%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'-0.82%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.278"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.33%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.2005"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'11.71%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.04492"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-2.38%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001214"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-1.35%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004636"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'4.87%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'0.16%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01611"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-5.91%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04335"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-8.98%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007403"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.43%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1302"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-3.72%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002007"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-14.70%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01034"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-7.79%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005895"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-1.60%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'0.16%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.220"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'170.60%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003008"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-11.34%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'0.16%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'0.16%"
$ws.Range("E50").Style = "Normal"
